$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing FY data (D:K) right to (E:L)
$ws.Range("D:D").Insert()

# Copy the number formats/styles from the (now shifted) column E into the
# newly inserted blank column D so each row picks up the same style it had
# before the insert (date style for header rows, numeric style for data rows).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore a sensible column width for the newly inserted column (it starts
# out at the generic default width after Insert).
$ws.Range("D1").ColumnWidth = $ws.Range("E1").ColumnWidth

# Populate the new column D with the FY2018 (period ending 2018-12-31) data.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 7014600
$ws.Range("D9").Value = 5369300
$ws.Range("D10").Value = 1645300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 17200
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 5946900
$ws.Range("D18").Value = 1067700
$ws.Range("D20").Value = 8000
$ws.Range("D21").Value = 1486600
$ws.Range("D22").Value = 105200
$ws.Range("D23").Value = 970500
$ws.Range("D24").Value = 234500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 736000
$ws.Range("D27").Value = 730300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 2000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -8000
$ws.Range("D33").Value = 732300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 732300

$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 361500
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 918600
$ws.Range("D44").Value = 795600
$ws.Range("D45").Value = 39400
$ws.Range("D46").Value = 2115100
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 3108600
$ws.Range("D49").Value = 1295500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 50500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 6569700
$ws.Range("D57").Value = 382200
$ws.Range("D58").Value = 1400
$ws.Range("D59").Value = 310000
$ws.Range("D60").Value = 693600
$ws.Range("D61").Value = 2501300
$ws.Range("D62").Value = 702400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3897300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 2315800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2672400
$ws.Range("D77").Value = 0

$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 732300
$ws.Range("D83").Value = 410900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1180100
$ws.Range("D91").Value = -555900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -608200
$ws.Range("D96").Value = -268100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -427300
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 144600

Write-Host "Done"
